# StoryCards.xlsx — "Bearbeitung StoryCards.xlsx nach falsch eingabe"
#
# Row 16 ("Überarbeitung GUI" / "Überarbeitung der GUI mit MIcrosoftBlend" / "12h")
# was entered by mistake, so it is deleted; rows 17-21 shift up to 16-20.
# The autofilter / named range / conditional formatting / data validation
# ranges that covered the table all shrink by one row accordingly, and the
# status colour-coding rules (fertig / in Arbeit / jungfräulich) get
# reordered (jungfräulich first, fertig last) while keeping the same
# effective colours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Remove the erroneously entered row.
$ws.Rows("16:16").Delete()

# 2. Re-apply the AutoFilter so its range shrinks from A7:M35 to A7:M34.
$ws.AutoFilterMode = $false
$ws.Range("A7:M34").AutoFilter()

# 3. Fix up the (now stale) _FilterDatabase defined name so it matches.
$wb.Names.Item("Tabelle1!_FilterDatabase").RefersTo = "=Tabelle1!`$A`$7:`$M`$34"

# 4. Reorder / recolour the status conditional-formatting rules on B8:B34.
#    Rule 1 currently tests "fertig" (green), rule 3 currently tests
#    "jungfräulich" (red). Swap what text/colour each rule position uses so
#    the rule for "jungfräulich" now comes first (priority 1) and "fertig"
#    comes last (priority 3) - net visual colours per status stay the same.
$rng = $ws.Range("B8:B34")
$fcs = $rng.FormatConditions

$fcFertig = $fcs.Item(1)
$fcJungfraeulich = $fcs.Item(3)

$fcJungfraeulich_color = $fcJungfraeulich.Interior.Color
$fcFertig_color = $fcFertig.Interior.Color

$fcFertig.Text = "jungfräulich"
$fcFertig.Formula1 = '=NOT(ISERROR(SEARCH("jungfräulich",B8)))'
$fcFertig.Priority = 1
$fcFertig.Interior.Color = $fcJungfraeulich_color

$fcJungfraeulich.Text = "fertig"
$fcJungfraeulich.Formula1 = '=NOT(ISERROR(SEARCH("fertig",B8)))'
$fcJungfraeulich.Priority = 3
$fcJungfraeulich.Interior.Color = $fcFertig_color

# Make sure the rules now cover B8:B34 (AppliesTo updates the whole group).
$fcs.Item(1).ModifyAppliesToRange($rng)

# 5. Selection moved to A21 (first empty row right after the shrunk table).
$ws.Activate()
$ws.Range("A21").Select()
